$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4348.7754
$ws.Range("I132").Value = 4326.372
$ws.Range("K132").Value = 12979.116
$ws.Range("M132").Value = -10449.116

$ws.Range("H133").Value = 84309.60000000001
$ws.Range("J133").Value = 84309.60000000001
$ws.Range("L133").Value = 84309.60000000001
$ws.Range("N133").Value = -94429.60000000001

$ws.Range("H138").Value = 315792.56
$ws.Range("I138").Value = 488670.16
$ws.Range("J138").Value = 4612.8667
$ws.Range("K138").Value = 1466010.48
$ws.Range("L138").Value = 13838.6001
$ws.Range("M138").Value = -1460870.48
$ws.Range("N138").Value = -24118.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H45").Value = 117568
$ws.Range("I45").Value = 187343.81
$ws.Range("J45").Value = 7920.2856
$ws.Range("K45").Value = 187343.81
$ws.Range("L45").Value = 7920.2856
$ws.Range("M45").Value = -186966.81
$ws.Range("N45").Value = -8674.285599999999

$ws.Range("H61").Value = 5467.3477
$ws.Range("I61").Value = 5484.0454
$ws.Range("J61").Value = 5100
$ws.Range("K61").Value = 5484.0454
$ws.Range("L61").Value = 5100
$ws.Range("M61").Value = -5272.0454
$ws.Range("N61").Value = -5524

$ws.Range("H102").Value = 12986.046
$ws.Range("I102").Value = 17453.615
$ws.Range("K102").Value = 17453.615
$ws.Range("M102").Value = -15831.615

$ws.Range("H110").Value = 2077.9333
$ws.Range("I110").Value = 914
$ws.Range("J110").Value = 3823.8333
$ws.Range("K110").Value = 914
$ws.Range("L110").Value = 3823.8333
$ws.Range("M110").Value = 1131
$ws.Range("N110").Value = -7913.8333

$ws.Range("H132").Value = 3710
$ws.Range("I132").Value = 2833.8696
$ws.Range("K132").Value = 8501.6088
$ws.Range("M132").Value = -5971.6088

$ws.Range("H136").Value = 5467.3477
$ws.Range("I136").Value = 5484.0454
$ws.Range("J136").Value = 5100
$ws.Range("K136").Value = 16452.1362
$ws.Range("L136").Value = 15300
$ws.Range("M136").Value = -13902.1362
$ws.Range("N136").Value = -20400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25727
$ws.Range("I82").Value = 7587
$ws.Range("J82").Value = 57472
$ws.Range("K82").Value = 7587
$ws.Range("L82").Value = 57472
$ws.Range("M82").Value = -7204
$ws.Range("N82").Value = -58238

$ws.Range("H85").Value = 25727
$ws.Range("I85").Value = 7587
$ws.Range("J85").Value = 57472
$ws.Range("K85").Value = 7587
$ws.Range("L85").Value = 57472
$ws.Range("M85").Value = -6261
$ws.Range("N85").Value = -60124

$ws.Range("H116").Value = 68475
$ws.Range("J116").Value = 68475
$ws.Range("L116").Value = 68475
$ws.Range("N116").Value = -77653

$ws.Range("H134").Value = 6489.7666
$ws.Range("I134").Value = 6560.4644
$ws.Range("K134").Value = 19681.3932
$ws.Range("M134").Value = -17146.3932

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 18002.143
$ws.Range("J86").Value = 19334.5
$ws.Range("L86").Value = 19334.5
$ws.Range("N86").Value = -21580.5

$ws.Range("H89").Value = 18002.143
$ws.Range("J89").Value = 19334.5
$ws.Range("L89").Value = 96672.5
$ws.Range("N89").Value = -107904.5

$ws.Range("H141").Value = 281180.94
$ws.Range("I141").Value = 56864
$ws.Range("K141").Value = 56864
$ws.Range("M141").Value = -51684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1340.4
$ws.Range("J34").Value = 1350
$ws.Range("L34").Value = 4050
$ws.Range("N34").Value = -4218

$ws.Range("H39").Value = 839.9
$ws.Range("J39").Value = 2999.5
$ws.Range("L39").Value = 8998.5
$ws.Range("N39").Value = -9586.5

$ws.Range("H55").Value = 6130.727
$ws.Range("I55").Value = 1645
$ws.Range("J55").Value = 7127.5557
$ws.Range("K55").Value = 4935
$ws.Range("L55").Value = 21382.6671
$ws.Range("M55").Value = -4758
$ws.Range("N55").Value = -21736.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3346944.5
$ws.Range("I14").Value = 3346944.5
$ws.Range("K14").Value = 3346944.5
$ws.Range("M14").Value = -3346776.5

$ws.Range("H70").Value = 8907.134
$ws.Range("I70").Value = 7843.2
$ws.Range("K70").Value = 7843.2
$ws.Range("M70").Value = -7573.2

$ws.Range("H73").Value = 8907.134
$ws.Range("I73").Value = 7843.2
$ws.Range("K73").Value = 7843.2
$ws.Range("M73").Value = -6907.2

$ws.Range("H80").Value = 16130.556
$ws.Range("I80").Value = 25999.75
$ws.Range("J80").Value = 8235.200000000001
$ws.Range("K80").Value = 25999.75
$ws.Range("L80").Value = 8235.200000000001
$ws.Range("M80").Value = -25001.75
$ws.Range("N80").Value = -10231.2

$ws.Range("H83").Value = 16130.556
$ws.Range("I83").Value = 25999.75
$ws.Range("J83").Value = 8235.200000000001
$ws.Range("K83").Value = 129998.75
$ws.Range("L83").Value = 41176
$ws.Range("M83").Value = -125006.75
$ws.Range("N83").Value = -51160

$ws.Range("H140").Value = 78212.37
$ws.Range("I140").Value = 40709
$ws.Range("K140").Value = 40709
$ws.Range("M140").Value = -35529

$ws.Range("H141").Value = 66279.55
$ws.Range("J141").Value = 66608.3
$ws.Range("L141").Value = 66608.3
$ws.Range("N141").Value = -76968.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13162.471
$ws.Range("I22").Value = 25594.5
$ws.Range("J22").Value = 2111.7778
$ws.Range("K22").Value = 25594.5
$ws.Range("L22").Value = 2111.7778
$ws.Range("M22").Value = -25299.5
$ws.Range("N22").Value = -2701.7778

$ws.Range("H27").Value = 13162.471
$ws.Range("I27").Value = 25594.5
$ws.Range("J27").Value = 2111.7778
$ws.Range("K27").Value = 25594.5
$ws.Range("L27").Value = 2111.7778
$ws.Range("M27").Value = -25487.5
$ws.Range("N27").Value = -2325.7778

$ws.Range("H46").Value = 4495.5835
$ws.Range("I46").Value = 983.3333
$ws.Range("J46").Value = 5666.3335
$ws.Range("K46").Value = 983.3333
$ws.Range("L46").Value = 5666.3335
$ws.Range("M46").Value = -795.3333
$ws.Range("N46").Value = -6042.3335

$ws.Range("H100").Value = 6077.143
$ws.Range("I100").Value = 2912.5
$ws.Range("J100").Value = 10296.667
$ws.Range("K100").Value = 2912.5
$ws.Range("L100").Value = 10296.667
$ws.Range("M100").Value = -2371.5
$ws.Range("N100").Value = -11378.667

$ws.Range("H132").Value = 576135.5600000001
$ws.Range("I132").Value = 1065709.9
$ws.Range("K132").Value = 3197129.7
$ws.Range("M132").Value = -3194599.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 79166.10000000001
$ws.Range("J141").Value = 71406.78
$ws.Range("L141").Value = 71406.78
$ws.Range("N141").Value = -81766.78

Write-Output "Applied all edits."
